$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after the current row 8 ("Financial conditions & employment, US.xlsx"),
# i.e. before the old row 9 ("U.S GDP and GDI.xlsx"), pushing everything down.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Bank credit and M2 U.S.xlsx"
$ws.Range("A10").Value = "Total construction spending in the U.S, nominal and real.xlsx"
$ws.Range("A11").Value = "Bitcoin price change is fueled by global monetary growth.xlsx"

# Old row 9 ("U.S GDP and GDI.xlsx") is now at row 12.
# Insert 3 more new rows after row 12, before the old row 10 ("US and Global Monetary Aggregates.xlsx"),
# which is now at row 13.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = "Bitcoin price change is fueled by global monetary growth (with forecast).xlsx"
$ws.Range("A14").Value = "Personal savings U.S (BEA).xlsx"
$ws.Range("A15").Value = "Real Incomes U.S (BEA).xlsx"
